# Add the eSNHBV1 endogenous snake hepatitis B virus reference row, and
# (re)apply the AutoFilter over the data range, matching the upstream
# "hepadnairusCompound aligner" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15 (pushes the existing rows 15-43 down to 16-44)
$ws.Rows("15:15").Insert()

# Populate the new reference row:
#   sequence-ID / name  : eSNHBV1-con
#   full_name           : Endogenous snake hepatitis B virus 1
#   genus               : Herpetohepadnavirus
#   host_sci_name       : Colubroidea
#   host_common_name    : snakes
$ws.Range("A15").Value2 = "eSNHBV1-con"
$ws.Range("B15").Value2 = "eSNHBV1-con"
$ws.Range("C15").Value2 = "Endogenous snake hepatitis B virus 1"
$ws.Range("D15").Value2 = "Herpetohepadnavirus"
$ws.Range("E15").Value2 = "Colubroidea"
$ws.Range("F15").Value2 = "snakes"

# Match formatting of the plain (unshaded) rows at the bottom of the table
# (row 44 post-insert == the former last data row, African cichlid) for the
# cells that should carry the "plain" style rather than the style inherited
# from the row pushed down during the insert.
$ws.Range("E44").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("F44").Copy()
$ws.Range("F15").PasteSpecial(-4122)

# Re-apply the table AutoFilter across the (now one row larger) data range.
$ws.Range("A2:F46").AutoFilter() | Out-Null

# Leave the selection where it ended up after the edit, as in the source file.
$ws.Range("C15").Select() | Out-Null
